# "Power Storage" sheet — set Exis Unit (column E) to 0 where it still had
# a nonzero placeholder, and raise MaxInvest (column S, the line-load / max
# invest cap) to 100 for rows 7-11.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Power Storage")

$ws.Range("E7").Value = 0
$ws.Range("E10").Value = 0

$ws.Range("S7:S11").Value = 100

$ws.Activate()
$ws.Range("S8:S11").Select()
